$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update record data (row 2: Sergio Montero -> Efren Moreno with new DNI code) ---
$ws.Range("A2").Value = "AUT_JF_QA_001"
$ws.Range("B2").Value = "Efrén"
$ws.Range("C2").Value = "Moreno"

# --- Row 3 gains a DNI code and accent fix on the name ---
$ws.Range("A3").Value = "AUT_TF_QA_001"
$ws.Range("B3").Value = "Raúl"

# --- Row 4 gains a DNI code and the first name is replaced ---
$ws.Range("A4").Value = "AUT_ TF_QA_002"
$ws.Range("B4").Value = "Héctor"

# --- Widen column A now that it holds the longer AUT_* codes ---
$ws.Columns.Item(1).ColumnWidth = 15.6

# --- Remove the stray formatted-but-empty A6 cell (B6 keeps its style) ---
$ws.Range("A6").Clear()

# --- Move the active selection to C6 ---
[void]$ws.Range("C6").Select()
